# Button method update (scroll)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: limit name / new limit values
$ws.Range("B2").Value = "QR Payments"
$ws.Range("C2").Value = "90000"

# Remove rows 3 and 4 (old "Funds Transfer to HBL Account" / "Funds Transfer to Other Banks Account" rows)
$ws.Rows("3:4").Delete()

# Restore selection to B8 as in the saved file
$ws.Range("B8").Select()
